# v1.2 - verified the notification test cases
$wb = $excel.ActiveWorkbook

# 1) Mark the three reviewed rows on the main review sheet as "Closed" in
#    the "Reviewer verification" column (J) now that the test cases have
#    been re-checked after the update.
$reviewSheet = $wb.Worksheets.Item("LH_TC_NOTIFICATION_REVIEWS")
$reviewSheet.Range("J2").Value = "Closed"
$reviewSheet.Range("J3").Value = "Closed"
$reviewSheet.Range("J4").Value = "Closed"
$reviewSheet.Range("J4").Select()

# 2) Log the new version in the "Version History" table.
$historySheet = $wb.Worksheets.Item("Version History")
$historySheet.Range("A4").Value = "v1.2"
$historySheet.Range("B4").Value = "Gehad Ashry"
$historySheet.Range("C4").Value = "Verified Test cases after update"
$historySheet.Range("D4").Formula = "=DATE(2025,4,21)"

# 3) Leave the "Version History" tab active/selected, as the author did
#    after adding the new version row.
$historySheet.Activate()
$historySheet.Range("G9").Select()
